# Updated cryptos list values (price / volume) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.525.06'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +3.72%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.501.96'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.40%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.42'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.30%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '169.80'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +5.36%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.501.05'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.30%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.588'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +6.56%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.33'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.64%  '

$ws.Range('E11').Value = '  +4.82%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.438'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +3.52%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.110.48'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.40%  '

$ws.Range('E14').Value = '  -0.44%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '28.35'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.90%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000178'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.26%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '66.574.45'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.73%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.496.46'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.26%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.34'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.36%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.08'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.86%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '390.55'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.71%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.97'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.93%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '73.04'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.34%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.04%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.536'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.56%  '

$ws.Range('E26').Value = '  +6.08%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.46'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +9.89%  '

$ws.Range('E28').Value = '  +1.98%  '

$ws.Range('E29').Value = '  +0.34%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.34'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +5.42%  '

$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.49'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +6.64%  '

$ws.Range('E32').Value = '  +2.62%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '23.59'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.22%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.39'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.70%  '

$ws.Range('E35').Value = '  +0.01%  '

$ws.Range('E36').Value = '  +7.43%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '162.33'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.72%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.883'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.29%  '

$ws.Range('E39').Value = '  +4.93%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.84'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.83%  '

$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.70'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +6.35%  '

$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '27.67'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.85%  '

$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0746'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.80%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '26.50'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.89%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.807.70'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.15%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '43.16'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.62%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0311'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.13%  '

$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.52'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.20%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '354.09'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +5.27%  '

$ws.Range('E50').Value = '  +3.31%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '33.75'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +12.40%  '
